$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "report generated" timestamp in M6 to the new dump's date/time.
$ws.Range("M6").Value = 41439.32492519676

# Clear out the sample data row (row 12) that was only a placeholder dump row.
$ws.Range("A12").Value = ""
$ws.Range("C12").Value = ""
$ws.Range("E12").Value = ""
$ws.Range("G12").Value = ""
$ws.Range("I12").Value = ""
$ws.Range("K12").Value = ""
$ws.Range("M12").Value = ""
$ws.Range("O12").Value = ""
